$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Subject" values for the two data rows
$ws.Range("B2").Value = "Immuno&hema"
$ws.Range("B3").Value = "Immuno&hema"

# Style the Subject column (B2:B3): centered text on a light-gray fill
$rng = $ws.Range("B2:B3")
$rng.Font.Name = "Calibri"
$rng.Font.Size = 11
$rng.Interior.Color = 15790320     # RGB(240, 240, 240) = 0xF0F0F0
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.VerticalAlignment = -4108     # xlCenter

# Mirror the selection state captured in the saved file
$ws.Range("B2:B3").Select()
